$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-17 23:50:32", 34, 31, 6, 7, 6, 8, 0.001, 0.01, 0.003, 100, 512, 10, 7, 0.9117647058823529),
    @("2023-12-18 18:37:39", 6, 5, 2, 0, 1, 0, 0.001, 0.01, 0.003, 100, 512, 10, 7, 0.8333333333333334),
    @("2023-12-18 18:56:47", 2, 2, 1, 0, 0, 1, 0.001, 0.01, 0.003, 100, 512, 10, 7, 1),
    @("2023-12-19 00:43:14", 15, 12, 3, 1, 3, 5, 0.001, 0.01, 0.003, 100, 512, 10, 7, 0.8),
    @("2023-12-19 01:16:00", 7, 5, 1, 0, 2, 2, 0.001, 0.01, 0.003, 100, 512, 10, 7, 0.7142857142857143),
    @("2023-12-21 07:21:25", 11, 7, 1, 0, 2, 2, 0.001, 0.01, 0.003, 100, 512, 10, 7, 0.6363636363636364),
    @("2023-12-22 22:48:45", 11, 8, 2, 3, 2, 1, 0.001, 0.01, 0.003, 100, 512, 10, 7, 0.7272727272727273)
)

$startRow = 13
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
